$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
# Fill order matches the source workbook's shared-string build order:
# columns A,B,D,H,M row-by-row for rows 2-7, then column C for all rows,
# then column J for all rows, then column R for all rows, with E/I (numeric)
# interspersed per row.

# Row 2: AHC30716 - Certificate III in Horticulture
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "CERTIFICATE III IN HORTICULTURE"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I2").Value = 8700
$ws.Range("M2").Value = "TAS"

# Row 3: AHC40416 - Certificate IV in Horticulture
$ws.Range("A3").Value = "AHC40416"
$ws.Range("B3").Value = "110598E"
$ws.Range("D3").Value = "CERTIFICATE IV IN HORTICULTURE"
$ws.Range("E3").Value = 52
$ws.Range("H3").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I3").Value = 8700
$ws.Range("M3").Value = "TAS"

# Row 4: AHC51422 - Diploma of Agribusiness Management
$ws.Range("A4").Value = "AHC51422"
$ws.Range("B4").Value = "110774E"
$ws.Range("D4").Value = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E4").Value = 52
$ws.Range("H4").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I4").Value = 9200
$ws.Range("M4").Value = "TAS"

# Row 5: Package - Certificate III + Certificate IV
$ws.Range("A5").Value = "AHC30716 / AHC40416"
$ws.Range("B5").Value = "110597F / 110598E"
$ws.Range("D5").Value = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I5").Value = 16700
$ws.Range("M5").Value = "TAS"

# Row 6: Package - Certificate III + Diploma
$ws.Range("A6").Value = "AHC30716 / AHC51422"
$ws.Range("B6").Value = "110597F / 110774E"
$ws.Range("D6").Value = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E6").Value = 104
$ws.Range("H6").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I6").Value = 17200
$ws.Range("M6").Value = "TAS"

# Row 7: Package - Certificate IV + Diploma
$ws.Range("A7").Value = "AHC40416 / AHC51422"
$ws.Range("B7").Value = "110598E / 110774E"
$ws.Range("D7").Value = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E7").Value = 104
$ws.Range("H7").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I7").Value = 17200
$ws.Range("M7").Value = "TAS"

# Column C (department) for all rows
$ws.Range("C2").Value = "HORTICULTURE"
$ws.Range("C3").Value = "HORTICULTURE"
$ws.Range("C4").Value = "MANAGEMENT"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("C6").Value = "PACKAGES"
$ws.Range("C7").Value = "PACKAGES"

# Column J (tuitionDetail) for all rows
$ws.Range("J2").Value = "8,500 tuition fee + 200 handling fee"
$ws.Range("J3").Value = "8,500 tuition fee + 200 handling fee"
$ws.Range("J4").Value = "9,000 tuition fee + 200 handling fee"
$ws.Range("J5").Value = "16,500 tuition fee + 200 handling fee"
$ws.Range("J6").Value = "17,000 tuition fee + 200 handling fee"
$ws.Range("J7").Value = "17,000 tuition fee + 200 handling fee"

# Column R (promotionValidity) for all rows
$ws.Range("R2").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R3").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R4").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R5").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R6").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R7").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Formatting ---------------------------------------------------------
# Apply cell formatting in the same order Excel originally created the styles,
# so the generated cellXfs indices line up with the source workbook:
#   index 2 = wrapText only
#   index 3 = number format "#,##0" only
#   index 4 = number format "#,##0" + wrapText

# durationDetail (H) -> wrapText only
$ws.Range("H2:H7").WrapText = $true

# tuition (I) -> number format only
$ws.Range("I2:I7").NumberFormat = "#,##0"

# tuitionDetail (J) -> number format + wrapText
$ws.Range("J2:J7").NumberFormat = "#,##0"
$ws.Range("J2:J7").WrapText = $true

# Package rows' vetCode/cricosCode/name columns (multi-line content) -> wrapText only
$ws.Range("A5:B7").WrapText = $true
$ws.Range("D5:D7").WrapText = $true

# Row heights for the new data rows
$ws.Range("2:7").RowHeight = 45

# Restore the cursor/selection position as left by the author
$ws.Range("S14").Select()
